# Update "想去人数" (want-to-go count) column F figures that were refreshed
# at the latest data scrape (gh-pages output generated at 456a3b4).
#
# Sheet "展览" (sheet1) and "全部类型" (sheet4) both list the same events,
# so the same F-column bumps apply to both, just shifted by one row on the
# second sheet (it carries one extra row before the matching tail row).

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  Old = 1067; New = 1071 },
    @{ Row = 3;  Old = 352;  New = 356  },
    @{ Row = 5;  Old = 8681; New = 8687 },
    @{ Row = 9;  Old = 275;  New = 278  },
    @{ Row = 11; Old = 7;    New = 8    },
    @{ Row = 12; Old = 3526; New = 3537 },
    @{ Row = 15; Old = 73;   New = 74   },
    @{ Row = 16; Old = 1124; New = 1132 },
    @{ Row = 20; Old = 193;  New = 195  },
    @{ Row = 21; Old = 2252; New = 2259 }
)

# Sheet "展览": last updated row is F22 (46 -> 47)
$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.New
}
$ws1.Cells.Item(22, 6).Value = 47

# Sheet "全部类型": last updated row is F23 (46 -> 47)
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.New
}
$ws4.Cells.Item(23, 6).Value = 47
